# Disaggregation of commodity Copper
#
# 1) Rename the commodity label "Copper ores and concentrates" -> "Copper"
#    (row 7, column C, on every yearly worksheet).
# 2) For every yearly worksheet, the Min/hist/Max values in columns D/E/F
#    for rows 5 (Neodymium), 7 (Copper) and 8 (Raw silicon) are cyclically
#    rotated one column to the right (D->E, E->F, F->D).

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- rename the commodity label ---
    $ws.Range("C7").Value2 = "Copper"

    # --- cyclic shift D->E->F->D for rows 5, 7 and 8 ---
    foreach ($r in 5, 7, 8) {
        $dAddr = "D" + $r
        $eAddr = "E" + $r
        $fAddr = "F" + $r

        $dVal = $ws.Range($dAddr).Value2
        $eVal = $ws.Range($eAddr).Value2
        $fVal = $ws.Range($fAddr).Value2

        $ws.Range($eAddr).Value2 = $dVal
        $ws.Range($fAddr).Value2 = $eVal
        $ws.Range($dAddr).Value2 = $fVal
    }
}
